$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2799475832522333
$ws.Range("C2").Value = 0.2090321810254712
$ws.Range("D2").Value = -0.118148624057875
$ws.Range("E2").Value = -0.05271246304120579
$ws.Range("G2").Value = -0.04489725391781101
$ws.Range("H2").Value = 0.170632811130217
$ws.Range("I2").Value = -0.2492238838232901
$ws.Range("J2").Value = 0.3717155353279262
$ws.Range("K2").Value = 0.0660253734085456
$ws.Range("M2").Value = 0.0721141320311092
$ws.Range("N2").Value = 0.2535374338888151
$ws.Range("O2").Value = 0.2614604786978406
$ws.Range("P2").Value = -0.2199382078615827
$ws.Range("Q2").Value = -0.06485291678922628
$ws.Range("S2").Value = -0.4227554336590582
$ws.Range("T2").Value = 0.2933123469905549
$ws.Range("U2").Value = 0.009921272601481622
$ws.Range("V2").Value = -0.1610876116140179

$ws.Range("B3").Value = 0.009003659529304139
$ws.Range("C3").Value = 0.05114354262493992
$ws.Range("D3").Value = 0.2703189106102171
$ws.Range("E3").Value = 0.6228550295960428
$ws.Range("G3").Value = 0.6752950248477201
$ws.Range("H3").Value = 0.1070925184168591
$ws.Range("I3").Value = 0.01813576185122696
$ws.Range("J3").Value = 0.0003477731161933515
$ws.Range("K3").Value = 0.5378828208791165
$ws.Range("M3").Value = 0.5010580359233224
$ws.Range("N3").Value = 0.01800535102541572
$ws.Range("O3").Value = 0.0147104444533948
$ws.Range("P3").Value = 0.04016705954746833
$ws.Range("Q3").Value = 0.5378757269451848
$ws.Range("S3").Value = 0.00008002855799176172
$ws.Range("T3").Value = 0.000001163841448708466
$ws.Range("U3").Value = 0.871211386875466
$ws.Range("V3").Value = 0.007995765619599986
